$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1
$ws.Range("D1").Value = 2.6
$ws.Range("E1").Value = 3.4

# Row 3
$ws.Range("C3").Value = 3.1
$ws.Range("D3").Value = 3.2
$ws.Range("E3").Value = 2.25

# Row 4
$ws.Range("C4").Value = 1.44
$ws.Range("D4").Value = 4.5
$ws.Range("E4").Value = 6

# Row 5
$ws.Range("C5").Value = 3.5
$ws.Range("D5").Value = 3.3
$ws.Range("E5").Value = 2

# Row 6
$ws.Range("C6").Value = 2.35
$ws.Range("E6").Value = 2.9

# Row 7
$ws.Range("C7").Value = 5.25
$ws.Range("D7").Value = 3.6
$ws.Range("E7").Value = 1.61

# Row 8
$ws.Range("C8").Value = 2.25
$ws.Range("D8").Value = 3.35
$ws.Range("E8").Value = 3

# Row 9
$ws.Range("C9").Value = 1.36
$ws.Range("D9").Value = 4.5
$ws.Range("E9").Value = 8

# Row 10
$ws.Range("C10").Value = 2.4
$ws.Range("D10").Value = 3.2
$ws.Range("E10").Value = 2.85

# Update selection to E10
$ws.Range("E10").Select()
